$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-17
# from 2023-09-19 (45188) to 2023-09-20 (45189)
$ws.Range("C2:C17").Value = 45189
